# Refresh the crypto price/volume table (GitHub Actions scheduled update).
# Price (column D) and 1h volume change % (column E) are refreshed per coin.
# D6 ("142.20") needs an explicit text format first; otherwise Excel's COM
# layer auto-coerces the numeric-looking string to a Number and the
# formatted display drops the trailing zero (142.2 instead of 142.20).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.897.28'
$ws.Range('E2').Value = '  -0.72%  '
$ws.Range('D3').Value = '2.409.76'
$ws.Range('E3').Value = '  -0.59%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '561.85'
$ws.Range('E5').Value = '  +1.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.20'
$ws.Range('E6').Value = '  -0.92%  '
$ws.Range('E8').Value = '  -0.79%  '
$ws.Range('E9').Value = '  -0.91%  '
$ws.Range('E10').Value = '  -1.83%  '
$ws.Range('E11').Value = '  -2.98%  '
$ws.Range('D12').Value = '0.349'
$ws.Range('E12').Value = '  -0.73%  '
$ws.Range('D13').Value = '25.49'
$ws.Range('E13').Value = '  -3.07%  '
$ws.Range('E14').Value = '  -1.77%  '
$ws.Range('D15').Value = '2.841.49'
$ws.Range('E15').Value = '  -0.74%  '
$ws.Range('D16').Value = '61.763.62'
$ws.Range('E16').Value = '  -0.79%  '
$ws.Range('D17').Value = '2.395.71'
$ws.Range('E17').Value = '  -1.19%  '
$ws.Range('D18').Value = '11.22'
$ws.Range('E18').Value = '  +1.22%  '
$ws.Range('D19').Value = '320.67'
$ws.Range('E19').Value = '  -1.21%  '
$ws.Range('D20').Value = '6.83'
$ws.Range('E20').Value = '  +1.29%  '
$ws.Range('E21').Value = '  -1.59%  '
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('D23').Value = '65.49'
$ws.Range('E23').Value = '  +1.00%  '
$ws.Range('E24').Value = '  -2.71%  '
$ws.Range('D25').Value = '8.67'
$ws.Range('E25').Value = '  -4.63%  '
$ws.Range('D26').Value = '563.77'
$ws.Range('E26').Value = '  -1.65%  '
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('E28').Value = '  -1.14%  '
$ws.Range('D29').Value = '0.0₃0932'
$ws.Range('E29').Value = '  -0.94%  '
$ws.Range('D30').Value = '8.17'
$ws.Range('E30').Value = '  -2.71%  '
$ws.Range('E31').Value = '  -4.68%  '
$ws.Range('E32').Value = '  -0.87%  '
$ws.Range('E33').Value = '  +0.09%  '
$ws.Range('D34').Value = '1.51'
$ws.Range('E34').Value = '  -4.30%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  -1.77%  '
$ws.Range('D37').Value = '5.44'
$ws.Range('E37').Value = '  -4.89%  '
$ws.Range('D38').Value = '152.34'
$ws.Range('E38').Value = '  +1.32%  '
$ws.Range('E39').Value = '  -1.47%  '
$ws.Range('E40').Value = '  -1.63%  '
$ws.Range('D41').Value = '1.78'
$ws.Range('E41').Value = '  -4.96%  '
$ws.Range('D43').Value = '147.93'
$ws.Range('E43').Value = '  -2.04%  '
$ws.Range('D44').Value = '2.24'
$ws.Range('E44').Value = '  -3.99%  '
$ws.Range('D45').Value = '3.59'
$ws.Range('E45').Value = '  -1.35%  '
$ws.Range('D46').Value = '0.0527'
$ws.Range('E46').Value = '  -3.05%  '
$ws.Range('D47').Value = '19.84'
$ws.Range('E47').Value = '  -2.81%  '
$ws.Range('D48').Value = '0.589'
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('D49').Value = '0.0917'
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('E50').Value = '  -1.57%  '
$ws.Range('D51').Value = '11.52'
$ws.Range('E51').Value = '  +0.30%  '
